# Update the "End Date" column (column E) values on Sheet1.
# The original scrape stamped every row with the scrape date (12-22-2022)
# in the End Date column. This corrects it to reflect each company's
# actual fiscal period end date, based on the ticker/section the row
# belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row ranges (inclusive) in column E and the corrected End Date value.
$ranges = @(
    @{Start=2;   End=222;  Value="12-31-2021"},
    @{Start=223; End=280;  Value="9-30-2022"},
    @{Start=281; End=339;  Value="12-31-2021"},
    @{Start=340; End=399;  Value="6-30-2022"},
    @{Start=400; End=456;  Value="1-31-2022"},
    @{Start=457; End=502;  Value="12-31-2021"}
)

foreach ($r in $ranges) {
    $addr = "E" + $r.Start + ":E" + $r.End
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $r.Value
    $rng.NumberFormat = "General"
}
